# Atualização automática de preços de eletricidade
# Updates row 2 of the spot price sheet with newly scraped values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (date serial number, advances one day)
$ws.Range("A2").Value = 45876

# Hourly prices (columns B..Z)
$ws.Range("B2").Value  = 116.95
$ws.Range("C2").Value  = 108.66
$ws.Range("D2").Value  = 105
$ws.Range("E2").Value  = 108.62
$ws.Range("F2").Value  = 107.96
$ws.Range("G2").Value  = 108.59
$ws.Range("H2").Value  = 114.7
$ws.Range("I2").Value  = 119.28
$ws.Range("J2").Value  = 119.22
$ws.Range("K2").Value  = 80.43000000000001
$ws.Range("L2").Value  = 50
$ws.Range("M2").Value  = 26.39
$ws.Range("N2").Value  = 26.39
$ws.Range("O2").Value  = 28.42
$ws.Range("P2").Value  = 27.2
$ws.Range("Q2").Value  = 26.39
$ws.Range("R2").Value  = 27.2
$ws.Range("S2").Value  = 48.97
$ws.Range("T2").Value  = 70
$ws.Range("U2").Value  = 101.53
$ws.Range("V2").Value  = 115.42
$ws.Range("W2").Value  = 145.01
$ws.Range("X2").Value  = 125.16
$ws.Range("Y2").Value  = 110
$ws.Range("Z2").Value  = 84.06

# Daily slot summary columns
$ws.Range("AA2").Value = "20h-24h"
$ws.Range("AB2").Value = 123.9
$ws.Range("AC2").Value = "20h-22h"
$ws.Range("AD2").Value = 130.22
$ws.Range("AE2").Value = "22h-24h"
$ws.Range("AF2").Value = 117.58
$ws.Range("AG2").Value = "9h-18h"
